$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM data for rows 2-10: Sending cluster x Target cluster matrix
# (ECs/FAPs/MuSCs) x (ECs/FAPs/MuSCs), ligand Efna5 -> receptor Epha2.
# NOTE: hashtable keys are case-insensitive, so the row-number key is named
# "rowNum" (not "r") to avoid colliding with the column "R" data key.
$rowData = @(
  @{ rowNum=2; A="ECs"; B="Efna5"; C="Epha2"; D="ECs"; E=1; F=0.3333333333333333; G=0.05800433333333333; H=0.174013; I=0.02087975181349295; J=0.02087975181349295; K=3; L=1; M=10.858287; N=32.574861; O=0.5084025289165609; P=0.508402528916561; Q=0.6298276985769999; R=5.668449287193; S=0.01061531862512997; T=0.01061531862512997 },
  @{ rowNum=3; A="ECs"; B="Efna5"; C="Epha2"; D="FAPs"; E=1; F=0.3333333333333333; G=0.05800433333333333; H=0.174013; I=0.02087975181349295; J=0.02087975181349295; K=2; L=0.6666666666666666; M=0.09477133333333332; N=0.284314; O=0.004437346842596906; P=0.004437346842596906; Q=0.00549714800911111; R=0.04947433208199999; S=[double]"9.265070078380999e-05"; T=[double]"9.265070078380999e-05" },
  @{ rowNum=4; A="ECs"; B="Efna5"; C="Epha2"; D="MuSCs"; E=1; F=0.3333333333333333; G=0.05800433333333333; H=0.174013; I=0.02087975181349295; J=0.02087975181349295; K=3; L=1; M=10.404599; N=31.213797; O=0.4871601242408422; P=0.4871601242408422; Q=0.6035118285956667; R=5.431606457361; S=0.01017178248757918; T=0.01017178248757918 },
  @{ rowNum=5; A="FAPs"; B="Efna5"; C="Epha2"; D="ECs"; E=3; F=1; G=1.666083666666667; H=4.998251; I=0.5997381826733804; J=0.5997381826733805; K=3; L=1; M=10.858287; N=32.574861; O=0.5084025289165609; P=0.508402528916561; Q=18.090814618679; R=162.817331568111; S=0.304908408758969; T=0.3049084087589691 },
  @{ rowNum=6; A="FAPs"; B="Efna5"; C="Epha2"; D="FAPs"; E=3; F=1; G=1.666083666666667; H=4.998251; I=0.5997381826733804; J=0.5997381826733805; K=2; L=0.6666666666666666; M=0.09477133333333332; N=0.284314; O=0.004437346842596906; P=0.004437346842596906; Q=0.1578969705348889; R=1.421072734814; S=0.002661246331270531; T=0.002661246331270532 },
  @{ rowNum=7; A="FAPs"; B="Efna5"; C="Epha2"; D="MuSCs"; E=3; F=1; G=1.666083666666667; H=4.998251; I=0.5997381826733804; J=0.5997381826733805; K=3; L=1; M=10.404599; N=31.213797; O=0.4871601242408422; P=0.4871601242408422; Q=17.33493245211633; R=156.014392069047; S=0.2921685275831409; T=0.292168527583141 },
  @{ rowNum=8; A="MuSCs"; B="Efna5"; C="Epha2"; D="ECs"; E=3; F=1; G=1.053930333333333; H=3.161791; I=0.3793820655131266; J=0.3793820655131266; K=3; L=1; M=10.858287; N=32.574861; O=0.5084025289165609; P=0.508402528916561; Q=11.443878037339; R=102.994902336051; S=0.192878801532462; T=0.192878801532462 },
  @{ rowNum=9; A="MuSCs"; B="Efna5"; C="Epha2"; D="FAPs"; E=3; F=1; G=1.053930333333333; H=3.161791; I=0.3793820655131266; J=0.3793820655131266; K=2; L=0.6666666666666666; M=0.09477133333333332; N=0.284314; O=0.004437346842596906; P=0.004437346842596906; Q=0.09988238293044444; R=0.8989414463739999; S=0.001683449810542565; T=0.001683449810542565 },
  @{ rowNum=10; A="MuSCs"; B="Efna5"; C="Epha2"; D="MuSCs"; E=3; F=1; G=1.053930333333333; H=3.161791; I=0.3793820655131266; J=0.3793820655131266; K=3; L=1; M=10.404599; N=31.213797; O=0.4871601242408422; P=0.4871601242408422; Q=10.96572249226967; R=98.691502430427; S=0.1848198141701221; T=0.1848198141701221 }
)

$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rowData) {
  for ($i = 0; $i -lt $colLetters.Length; $i++) {
    $colNum = $i + 1
    $letter = $colLetters[$i]
    $ws.Cells.Item($row.rowNum, $colNum).Value = $row[$letter]
  }
}

Write-Output "Updated rows 2-10 with new TPM values"